# Applies the "wall merging" requirements-doc colour-coding pass:
#  - a set of bullet paragraphs in the General/Scene/Settings sections get a
#    font colour applied (orange FF7F00 or green 00A800) to flag their status
#  - the "File | Edit | Module | ..." menu-structure line gets its 15 runs
#    collapsed back into a single run while keeping its original colour.
#
# Word colour values are 0x00BBGGRR decimal, i.e.:
#   FF7F00 (orange) -> 32767
#   00A800 (green)  -> 43008

$d = $word.ActiveDocument

$orange = 32767   # FF7F00
$green  = 43008   # 00A800

$colorMap = @{
    2  = $orange
    3  = $green
    4  = $orange
    5  = $orange
    6  = $green
    16 = $green
    17 = $orange
    22 = $green
    23 = $orange
    24 = $orange
    26 = $orange
    27 = $orange
    29 = $orange
    30 = $green
    31 = $green
    32 = $orange
    33 = $orange
    34 = $orange
    52 = $green
}

foreach ($idx in $colorMap.Keys) {
    $p = $d.Paragraphs($idx)
    $p.Range.Font.Color = $colorMap[$idx]
}

# Collapse the "File | Edit | ..." paragraph's many runs into a single run,
# preserving its original (dark) colour, without touching the paragraph mark.
$menuPara = $d.Paragraphs(55)
$menuRange = $menuPara.Range
$menuRange.End = $menuRange.End - 1
$menuColor = $menuRange.Font.Color
$menuStart = $menuRange.Start
$menuText = $menuRange.Text
$menuRange.Text = ""
$menuRange.InsertAfter($menuText)
$newMenuRange = $d.Range($menuStart, $menuRange.End)
$newMenuRange.Font.Color = $menuColor
